# Fruta / hortaliza, semanal
# Insert one new weekly price-record row at row 475 (pushing the existing
# rows 475-495 down to 476-496), matching the new row of data that was
# added to the source dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 475..495 down by one row, creating a blank row 475.
$ws.Rows.Item(475).Insert()

# Fill the new row 475 with the new record's data.
$ws.Range("A475").Value2 = 5
$ws.Range("B475").Value2 = "Macroferia Regional de Talca"
$ws.Range("C475").Value2 = "Maule"
$ws.Range("D475").Value2 = 45267
$ws.Range("E475").Value2 = 7
$ws.Range("F475").Value2 = "Fruta"
$ws.Range("G475").Value2 = 100108
$ws.Range("H475").Value2 = "Tropicales y subtropicales"
$ws.Range("I475").Value2 = 100108005
$ws.Range("J475").Value2 = "Piña"
$ws.Range("K475").Value2 = "Caramelo"
$ws.Range("L475").Value2 = "Tercera"
$ws.Range("M475").Value2 = 200
$ws.Range("N475").Value2 = 25000
$ws.Range("O475").Value2 = 25000
$ws.Range("P475").Value2 = 25000
$ws.Range("Q475").Value2 = "$/caja 16 unidades"
$ws.Range("R475").Value2 = "Ecuador"
$ws.Range("S475").Value2 = 1562
$ws.Range("T475").Value2 = 16
